$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: account holder name / card number
$ws.Range("C2").Value = "Hartmut"
# B3 holds a long all-digit card number that must stay TEXT (like the original
# inlineStr cell) rather than being auto-coerced to a number. A plain
# Value="2570314725427075" assignment would store it as a numeric cell, so we
# prefix with an apostrophe (Excel's "store as text" quote-prefix convention)
# to force a text cell, then re-apply row 2's formatting on top so the cell
# keeps its original look instead of the quote-prefix style.
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance statement
$ws.Range("D5").Value = "KONTOSTAND AM 04.09.2024"

# Row 6 (existing transaction line, text updated)
$ws.Range("B6").Value = "08.09."
$ws.Range("C6").Value = "09.09."
$ws.Range("D6").Value = "PAYPAL RRCXZG"
$ws.Range("E6").Value = "72,43-"

# Row 7 (existing transaction line, text updated)
$ws.Range("B7").Value = "12.09."
$ws.Range("C7").Value = "13.09."
$ws.Range("D7").Value = "KARTENZ./12.09 EDEKA RO"
$ws.Range("E7").Value = "33,03-"

# Row 8 (existing transaction line, text updated)
$ws.Range("B8").Value = "15.09."
$ws.Range("C8").Value = "16.09."
$ws.Range("D8").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 97801438"
$ws.Range("E8").Value = "86,00-"

# Rows 9-11 were previously blank placeholder rows; they now carry new
# transaction entries. Copy the formatting from the amount column of row 8
# (style index 17 in the original workbook) onto E9:E11 before writing values,
# so the "BETRAG" column keeps the same right-aligned look as the other
# transaction rows instead of its old placeholder style.
$ws.Range("E8").Copy() | Out-Null
$ws.Range("E9:E11").PasteSpecial(-4122) | Out-Null

# Row 9
$ws.Range("B9").Value = "18.09."
$ws.Range("C9").Value = "19.09."
$ws.Range("D9").Value = "KARTENZ./18.09 LIDL RO"
$ws.Range("E9").Value = "19,74-"

# Row 10
$ws.Range("B10").Value = "20.09."
$ws.Range("C10").Value = "21.09."
$ws.Range("D10").Value = "KARTENZ./20.09 LIDL RO"
$ws.Range("E10").Value = "129,92-"

# Row 11
$ws.Range("B11").Value = "22.09."
$ws.Range("C11").Value = "23.09."
$ws.Range("D11").Value = "KARTENZ./22.09 REWE RO"
$ws.Range("E11").Value = "21,64-"

# Row 12: closing balance statement
$ws.Range("D12").Value = "KONTOSTAND AM 25.09.2024"
$ws.Range("E12").Value = "362,76-"

# Row 13: next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 02.10.2024"
